# Add a new "edge" row to the browser/url table, mirroring the formatting
# of the row above it (row 3), then set values and fix the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 3's formatting (borders, alignment, font) down into row 4 so the
# new cells pick up the same styles as the existing rows instead of
# defaulting to unstyled cells.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data row.
$ws.Range("A4").Value = "edge"
$ws.Range("B4").Value = "https://demo.nopcommerce.com/"

# Match the saved selection state.
$ws.Range("B8").Select()
